$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '58.766.48'
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = '  -2.98%  '

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '2.723.14'
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = '  -6.21%  '

# Row 4
$ws.Range("E4").Value = '  -0.03%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '504.91'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -4.40%  '

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.34'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  -0.97%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = '  -0.29%  '

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.532'
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = '  -3.89%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '2.739.76'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  -5.71%  '

# Row 10
$ws.Range("E10").Value = '  -2.18%  '

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.09'
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = '  +2.51%  '

# Row 12
$ws.Range("E12").Value = '  -2.74%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '3.196.55'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -6.21%  '

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '58.761.06'
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = '  -2.96%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.72'
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = '  -3.90%  '

# Row 17
$ws.Range("E17").Value = '  -4.08%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '2.727.21'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -6.23%  '

# Row 19
$ws.Range("E19").Value = '  -5.35%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '11.02'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  -5.65%  '

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '343.04'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -5.69%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.27'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -4.80%  '

# Row 23
$ws.Range("E23").Value = '  -0.20%  '

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '5.64'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +0.02%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '63.39'
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = '  -1.15%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.175'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.38%  '

# Row 27
$ws.Range("E27").Value = '  -5.22%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.995'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.35%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.53'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  -3.42%  '

# Row 30
$ws.Range("E30").Value = '  -3.08%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.998'
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = '  -0.09%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '19.29'
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = '  -1.44%  '

# Row 33
$ws.Range("E33").Value = '  -4.14%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.72'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.95%  '

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.23'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  -3.13%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.45'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  -2.30%  '

# Row 37
$ws.Range("E37").Value = '  -5.73%  '

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.14'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -5.52%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '36.20'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -4.56%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.60'
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = '  -1.90%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '1.40'
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = '  -6.74%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '2.191.24'
$ws.Range("D42").ClearFormats()

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.0564'
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = '  -1.75%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.996'
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = '  -0.06%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.602'
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = '  -6.35%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '19.13'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -7.18%  '

# Row 47
$ws.Range("B47").Value = 'RenderToken'
$ws.Range("C47").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '4.80'
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = '  -4.00%  '

# Row 48
$ws.Range("B48").Value = 'WhiteBITCoin'
$ws.Range("C48").Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '10.37'
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = '  +0.15%  '

# Row 49
$ws.Range("E49").Value = '  -2.77%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0890'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -4.52%  '

# Row 51
$ws.Range("E51").Value = '  -1.66%  '
